# "pdf viewer con HTML5" — add the next sprint's backlog rows (17-21) to the
# estimation sheet: two 0.5-day tasks, a user-story row (merged, styled like
# the existing story rows but in a lighter grey), then two more tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: "Crear catálogo de usuarios" — 0.5 days
$ws.Range("B17").Value = 0.5
$ws.Range("C17").Value = "Crear catálogo de usuarios"

# Row 18: "Agregar busqueda por nombre, nick y email" — 0.5 days
$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = "Agregar busqueda por nombre, nick y email"

# Row 19: the user-story cell, merged B19:I19, formatted like the other
# story row (B13) but with a lighter grey font afterwards.
$ws.Range("B19").Value = "Como dueño del negocio, quiero ver los reportes de ventas del día, de días anteriores y comparativas, en tiempo real y en ambiente web, para tomar decisiones oportunas sin importar donde me encuentre."
$ws.Range("B13:I13").Copy() | Out-Null
$ws.Range("B19:I19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19:I19").Font.Color = 8421504
$ws.Range("B19:I19").Merge()
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Crear página reportes" — 0.5 days
$ws.Range("B20").Value = 0.5
$ws.Range("C20").Value = "Crear página reportes"

# Row 21: "crear pdf de prueba(...)" — 1 day
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "crear pdf de prueba(incluyendo filtros y visualizacion en un tab del sistema)"

# Move the selection to reflect where the user ended up after typing the new rows.
$ws.Range("B22").Select() | Out-Null

Write-Host "done"
